$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 label changes from "Service Life" to "service_lives"
$ws.Range("B1").Value = "service_lives"

# Update the active selection to S12 as in the target sheet view
$ws.Range("S12").Select()
